# Update cached market-price / profit figures across sheets (scheduled runner refresh)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 598.0476
$ws.Range("I12").Value = 558.8889
$ws.Range("K12").Value = 558.8889
$ws.Range("M12").Value = -388.8889

$ws.Range("H17").Value = 2789.6
$ws.Range("J17").Value = 3237.25
$ws.Range("L17").Value = 9711.75
$ws.Range("N17").Value = -10047.75

$ws.Range("H33").Value = 392.27274
$ws.Range("J33").Value = 834
$ws.Range("L33").Value = 834
$ws.Range("N33").Value = -1292

$ws.Range("H43").Value = 8298.956
$ws.Range("I43").Value = 8999
$ws.Range("K43").Value = 8999
$ws.Range("M43").Value = -8930

$ws.Range("H112").Value = 2094.8838
$ws.Range("J112").Value = 2246.0789
$ws.Range("L112").Value = 6738.236699999999
$ws.Range("N112").Value = -8954.236699999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 749.5
$ws.Range("I17").Value = 699
$ws.Range("J17").Value = 800
$ws.Range("K17").Value = 699
$ws.Range("L17").Value = 800
$ws.Range("M17").Value = -526
$ws.Range("N17").Value = -1146

$ws.Range("H74").Value = 4624.7827
$ws.Range("I74").Value = 4841.875
$ws.Range("J74").Value = 4128.5713
$ws.Range("K74").Value = 4841.875
$ws.Range("L74").Value = 4128.5713
$ws.Range("M74").Value = -3967.875
$ws.Range("N74").Value = -5876.5713

$ws.Range("H77").Value = 4624.7827
$ws.Range("I77").Value = 4841.875
$ws.Range("J77").Value = 4128.5713
$ws.Range("K77").Value = 24209.375
$ws.Range("L77").Value = 20642.8565
$ws.Range("M77").Value = -19841.375
$ws.Range("N77").Value = -29378.8565

$ws.Range("H132").Value = 780877.2
$ws.Range("I132").Value = 984071.4399999999
$ws.Range("J132").Value = 110336.3
$ws.Range("K132").Value = 2952214.32
$ws.Range("L132").Value = 331008.9
$ws.Range("M132").Value = -2949684.32
$ws.Range("N132").Value = -336068.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5838.75
$ws.Range("I86").Value = 1896.5714
$ws.Range("K86").Value = 1896.5714
$ws.Range("M86").Value = -773.5714

$ws.Range("H89").Value = 5838.75
$ws.Range("I89").Value = 1896.5714
$ws.Range("K89").Value = 9482.857
$ws.Range("M89").Value = -3866.857

$ws.Range("H134").Value = 1226657.4
$ws.Range("I134").Value = 1476263.4
$ws.Range("K134").Value = 4428790.199999999
$ws.Range("M134").Value = -4426255.199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 21413.857
$ws.Range("J3").Value = 19949.5
$ws.Range("L3").Value = 19949.5
$ws.Range("N3").Value = -20175.5

$ws.Range("H4").Value = 5454.364
$ws.Range("I4").Value = 4888.6665
$ws.Range("J4").Value = 8000
$ws.Range("K4").Value = 4888.6665
$ws.Range("L4").Value = 8000
$ws.Range("M4").Value = -4776.6665
$ws.Range("N4").Value = -8224

$ws.Range("H5").Value = 791.3333
$ws.Range("I5").Value = 337
$ws.Range("J5").Value = 1700
$ws.Range("K5").Value = 337
$ws.Range("L5").Value = 1700
$ws.Range("M5").Value = -225
$ws.Range("N5").Value = -1924

$ws.Range("H10").Value = 912
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()

$ws.Range("H22").Value = 1116631.2
$ws.Range("I22").Value = 1623723.8
$ws.Range("J22").Value = 1027.8
$ws.Range("K22").Value = 1623723.8
$ws.Range("L22").Value = 1027.8
$ws.Range("M22").Value = -1623373.8
$ws.Range("N22").Value = -1727.8

$ws.Range("H31").Value = 6008.3
$ws.Range("I31").Value = 2519.2
$ws.Range("K31").Value = 2519.2
$ws.Range("M31").Value = -2224.2

$ws.Range("H34").Value = 6008.3
$ws.Range("I34").Value = 2519.2
$ws.Range("K34").Value = 2519.2
$ws.Range("M34").Value = -2317.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 4757.12
$ws.Range("I122").Value = 713.8333
$ws.Range("J122").Value = 6033.9473
$ws.Range("K122").Value = 6424.4997
$ws.Range("L122").Value = 54305.5257
$ws.Range("M122").Value = -3974.4997
$ws.Range("N122").Value = -59205.5257

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 63.75
$ws.Range("J2").Value = 84.166664
$ws.Range("L2").Value = 84.166664
$ws.Range("N2").Value = -310.166664

$ws.Range("H70").Value = 7055.4136
$ws.Range("I70").Value = 7230.8
$ws.Range("J70").Value = 6867.5
$ws.Range("K70").Value = 7230.8
$ws.Range("L70").Value = 6867.5
$ws.Range("M70").Value = -6960.8
$ws.Range("N70").Value = -7407.5

$ws.Range("H73").Value = 7055.4136
$ws.Range("I73").Value = 7230.8
$ws.Range("J73").Value = 6867.5
$ws.Range("K73").Value = 7230.8
$ws.Range("L73").Value = 6867.5
$ws.Range("M73").Value = -6294.8
$ws.Range("N73").Value = -8739.5

$ws.Range("H80").Value = 23813194
$ws.Range("I80").Value = 58825850
$ws.Range("J80").Value = 4588.36
$ws.Range("K80").Value = 58825850
$ws.Range("L80").Value = 4588.36
$ws.Range("M80").Value = -58824852
$ws.Range("N80").Value = -6584.36

$ws.Range("H83").Value = 23813194
$ws.Range("I83").Value = 58825850
$ws.Range("J83").Value = 4588.36
$ws.Range("K83").Value = 294129250
$ws.Range("L83").Value = 22941.8
$ws.Range("M83").Value = -294124258
$ws.Range("N83").Value = -32925.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3105.8125
$ws.Range("I16").Value = 362
$ws.Range("J16").Value = 11337.25
$ws.Range("K16").Value = 362
$ws.Range("L16").Value = 11337.25
$ws.Range("M16").Value = -192
$ws.Range("N16").Value = -11677.25

$ws.Range("H61").Value = 5444.543
$ws.Range("I61").Value = 4567.087
$ws.Range("J61").Value = 7126.3335
$ws.Range("K61").Value = 4567.087
$ws.Range("L61").Value = 7126.3335
$ws.Range("M61").Value = -4365.087
$ws.Range("N61").Value = -7530.3335

$ws.Range("H113").Value = 5444.543
$ws.Range("I113").Value = 4567.087
$ws.Range("J113").Value = 7126.3335
$ws.Range("K113").Value = 4567.087
$ws.Range("L113").Value = 7126.3335
$ws.Range("M113").Value = -2397.087
$ws.Range("N113").Value = -11466.3335

$ws.Range("H132").Value = 4901.75
$ws.Range("I132").Value = 4901.75
$ws.Range("K132").Value = 14705.25
$ws.Range("M132").Value = -12175.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 17025
$ws.Range("J25").Value = 17025
$ws.Range("L25").Value = 17025
$ws.Range("N25").Value = -17611

$ws.Range("H33").Value = 24222
$ws.Range("I33").Value = 24222
$ws.Range("K33").Value = 24222
$ws.Range("M33").Value = -23972

$ws.Range("H36").Value = 24222
$ws.Range("I36").Value = 24222
$ws.Range("K36").Value = 24222
$ws.Range("M36").Value = -23972

$ws.Range("H126").Value = 6364.7144
$ws.Range("I126").Value = 3760.6
$ws.Range("K126").Value = 11281.8
$ws.Range("M126").Value = -8811.799999999999

$ws.Range("H132").Value = 9064
$ws.Range("I132").Value = 6816.6
$ws.Range("J132").Value = 17491.75
$ws.Range("K132").Value = 20449.8
$ws.Range("L132").Value = 52475.25
$ws.Range("M132").Value = -17919.8
$ws.Range("N132").Value = -57535.25
